$wb = $excel.ActiveWorkbook

# --- Sheet "표준식" : the sheet that actually receives the content updates ---
$ws = $wb.Worksheets.Item("표준식")

# Row 2 ("UI구현"): fill in description, assigned/finished dates, result and problem
$ws.Range("B2").Value = "기존에 command line으로 실행하던 프로그램을 웹기반 gui로 재구성"
$ws.Range("C2").Value = "05/22/2019"
$ws.Range("D2").Value = "05/29/2019"
$ws.Range("E2").Value = "웹 기반 GUI 구현"
$ws.Range("F2").Value = "DB 및 서버와 연결 필요"

# Row 3: replace "code refactoring" task with the new DB/server-connection task
$ws.Range("A3").Value = "DB연결 및 서버와 연결"
$ws.Range("B3").VerticalAlignment = -4107   # xlVAlignBottom -> drop the "top" alignment, keep wrap
$ws.Range("C3").Value = "05/31/2019"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""

# Row 4 ("게시판 분리"): give it a start date
$ws.Range("C4").Value = "05/31/2019"

# Row 5 ("게시판 분리", phonetic variant): give it a start date
$ws.Range("C5").Value = "05/31/2019"

# Row heights following the content changes above
$ws.Rows.Item(2).RowHeight = 34.5
$ws.Rows.Item(3).RowHeight = 42.75
$ws.Rows.Item(4).RowHeight = 51.75

# Column widths to fit the new text
$ws.Columns.Item(2).ColumnWidth = 42.57142857142857
$ws.Columns.Item(5).ColumnWidth = 14.857142857142858

# --- Selections on the other sheets / final active sheet+cell ---
$ws1 = $wb.Worksheets.Item("윤다은")
$ws1.Range("B2").Select() | Out-Null

$ws2 = $wb.Worksheets.Item("문준범")
$ws2.Range("D20").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("C4").Select() | Out-Null
